$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC!row100 (G=19906)
$ws_ALC.Range("H100").Value = 12822305
$ws_ALC.Range("I100").Value = 12822305
$ws_ALC.Range("J100").Value = 0
$ws_ALC.Range("K100").Value = 12822305
$ws_ALC.Range("L100").Value = 0
$ws_ALC.Range("M100").ClearContents()
$ws_ALC.Range("N100").Value = -12821764

# ALC!row116 (G=27778)
$ws_ALC.Range("H116").Value = 5196.8623
$ws_ALC.Range("I116").Value = 6492.316
$ws_ALC.Range("J116").Value = 2735.5
$ws_ALC.Range("K116").Value = 6492.316
$ws_ALC.Range("L116").Value = 2735.5
$ws_ALC.Range("M116").Value = -3050.316
$ws_ALC.Range("N116").Value = -9619.5

# ALC!row137 (G=44013)
$ws_ALC.Range("H137").Value = 1674.125
$ws_ALC.Range("I137").Value = 1282.909
$ws_ALC.Range("K137").Value = 3848.727
$ws_ALC.Range("M137").Value = -1298.727

# ARM!row32 (G=44147)
$ws_ARM.Range("H32").Value = 5540.311
$ws_ARM.Range("I32").Value = 4169.684
$ws_ARM.Range("J32").Value = 10135.941
$ws_ARM.Range("K32").Value = 4169.684
$ws_ARM.Range("L32").Value = 10135.941
$ws_ARM.Range("M32").Value = -3882.684
$ws_ARM.Range("N32").Value = -10709.941

# BSM!row34 (G=2421)
$ws_BSM.Range("H34").Value = 5026
$ws_BSM.Range("J34").Value = 5026
$ws_BSM.Range("L34").Value = 5026
$ws_BSM.Range("N34").Value = -5254

# BSM!row86 (G=12526)
$ws_BSM.Range("H86").Value = 41668416
$ws_BSM.Range("I86").Value = 47620764
$ws_BSM.Range("J86").Value = 2000
$ws_BSM.Range("K86").Value = 47620764
$ws_BSM.Range("L86").Value = 2000
$ws_BSM.Range("M86").Value = -47619641
$ws_BSM.Range("N86").Value = -4246

# BSM!row89 (G=12526)
$ws_BSM.Range("H89").Value = 41668416
$ws_BSM.Range("I89").Value = 47620764
$ws_BSM.Range("J89").Value = 2000
$ws_BSM.Range("K89").Value = 238103820
$ws_BSM.Range("L89").Value = 10000
$ws_BSM.Range("M89").Value = -238098204
$ws_BSM.Range("N89").Value = -21232

# BSM!row105 (G=19947)
$ws_BSM.Range("H105").Value = 38463572
$ws_BSM.Range("I105").Value = 62501604
$ws_BSM.Range("J105").Value = 2720
$ws_BSM.Range("K105").Value = 62501604
$ws_BSM.Range("L105").Value = 2720
$ws_BSM.Range("M105").Value = -62499857
$ws_BSM.Range("N105").Value = -6214

# BSM!row134 (G=43998)
$ws_BSM.Range("H134").Value = 4410.2144
$ws_BSM.Range("I134").Value = 6296.174
$ws_BSM.Range("J134").Value = 2127.2104
$ws_BSM.Range("K134").Value = 18888.522
$ws_BSM.Range("L134").Value = 6381.6312
$ws_BSM.Range("M134").Value = -16353.522
$ws_BSM.Range("N134").Value = -11451.6312

# CRP!row31 (G=44023)
$ws_CRP.Range("H31").Value = 12827225
$ws_CRP.Range("I31").Value = 2672.3125
$ws_CRP.Range("J31").Value = 21748654
$ws_CRP.Range("K31").Value = 2672.3125
$ws_CRP.Range("L31").Value = 21748654
$ws_CRP.Range("M31").Value = -2377.3125
$ws_CRP.Range("N31").Value = -21749244

# CRP!row34 (G=44023)
$ws_CRP.Range("H34").Value = 12827225
$ws_CRP.Range("I34").Value = 2672.3125
$ws_CRP.Range("J34").Value = 21748654
$ws_CRP.Range("K34").Value = 2672.3125
$ws_CRP.Range("L34").Value = 21748654
$ws_CRP.Range("M34").Value = -2470.3125
$ws_CRP.Range("N34").Value = -21749058

# CRP!row95 (G=18192)
$ws_CRP.Range("H95").Value = 35000
$ws_CRP.Range("J95").Value = 35000
$ws_CRP.Range("L95").Value = 35000
$ws_CRP.Range("N95").Value = -40492

# CRP!row99 (G=36198)
$ws_CRP.Range("H99").Value = 7362142
$ws_CRP.Range("I99").Value = 10041.2
$ws_CRP.Range("J99").Value = 17865142
$ws_CRP.Range("K99").Value = 10041.2
$ws_CRP.Range("L99").Value = 17865142
$ws_CRP.Range("M99").Value = -8543.200000000001
$ws_CRP.Range("N99").Value = -17868138

# CRP!row106 (G=18661)
$ws_CRP.Range("H106").Value = 41333.332
$ws_CRP.Range("J106").Value = 41333.332
$ws_CRP.Range("L106").Value = 41333.332
$ws_CRP.Range("N106").Value = -43857.332

# CRP!row126 (G=36198)
$ws_CRP.Range("H126").Value = 7362142
$ws_CRP.Range("I126").Value = 10041.2
$ws_CRP.Range("J126").Value = 17865142
$ws_CRP.Range("K126").Value = 30123.6
$ws_CRP.Range("L126").Value = 53595426
$ws_CRP.Range("M126").Value = -27653.6
$ws_CRP.Range("N126").Value = -53600366

# CRP!row134 (G=44020)
$ws_CRP.Range("H134").Value = 8774500
$ws_CRP.Range("I134").Value = 12348603
$ws_CRP.Range("J134").Value = 1702
$ws_CRP.Range("K134").Value = 37045809
$ws_CRP.Range("L134").Value = 5106
$ws_CRP.Range("M134").Value = -37043274
$ws_CRP.Range("N134").Value = -10176

# CUL!row68 (G=12895)
$ws_CUL.Range("H68").Value = 1240.5652
$ws_CUL.Range("I68").Value = 510.44446
$ws_CUL.Range("J68").Value = 1709.9286
$ws_CUL.Range("K68").Value = 1531.33338
$ws_CUL.Range("L68").Value = 5129.7858
$ws_CUL.Range("M68").Value = -720.33338
$ws_CUL.Range("N68").Value = -6751.7858

# CUL!row71 (G=12895)
$ws_CUL.Range("H71").Value = 1240.5652
$ws_CUL.Range("I71").Value = 510.44446
$ws_CUL.Range("J71").Value = 1709.9286
$ws_CUL.Range("K71").Value = 4594.00014
$ws_CUL.Range("L71").Value = 15389.3574
$ws_CUL.Range("M71").Value = -538.0001400000001
$ws_CUL.Range("N71").Value = -23501.3574

# CUL!row75 (G=12863)
$ws_CUL.Range("H75").Value = 47620372
$ws_CUL.Range("J75").Value = 47620372
$ws_CUL.Range("L75").Value = 142861116
$ws_CUL.Range("N75").Value = -142863112

# CUL!row78 (G=12863)
$ws_CUL.Range("H78").Value = 47620372
$ws_CUL.Range("J78").Value = 47620372
$ws_CUL.Range("L78").Value = 428583348
$ws_CUL.Range("N78").Value = -428593332

# CUL!row136 (G=44093)
$ws_CUL.Range("H136").Value = 28755
$ws_CUL.Range("I136").Value = 100030
$ws_CUL.Range("J136").Value = 4996.6665
$ws_CUL.Range("K136").Value = 300090
$ws_CUL.Range("L136").Value = 14989.9995
$ws_CUL.Range("M136").Value = -294990
$ws_CUL.Range("N136").Value = -25189.9995

# GSM!row69 (G=11891)
$ws_GSM.Range("H69").Value = 49600.5
$ws_GSM.Range("J69").Value = 49600.5
$ws_GSM.Range("L69").Value = 49600.5
$ws_GSM.Range("N69").Value = -51098.5

# GSM!row72 (G=11891)
$ws_GSM.Range("H72").Value = 49600.5
$ws_GSM.Range("J72").Value = 49600.5
$ws_GSM.Range("L72").Value = 148801.5
$ws_GSM.Range("N72").Value = -156289.5

# GSM!row80 (G=12521)
$ws_GSM.Range("H80").Value = 2647.6924
$ws_GSM.Range("I80").Value = 2482.5
$ws_GSM.Range("J80").Value = 2789.2856
$ws_GSM.Range("K80").Value = 2482.5
$ws_GSM.Range("L80").Value = 2789.2856
$ws_GSM.Range("M80").Value = -1484.5
$ws_GSM.Range("N80").Value = -4785.2856

# GSM!row83 (G=12521)
$ws_GSM.Range("H83").Value = 2647.6924
$ws_GSM.Range("I83").Value = 2482.5
$ws_GSM.Range("J83").Value = 2789.2856
$ws_GSM.Range("K83").Value = 12412.5
$ws_GSM.Range("L83").Value = 13946.428
$ws_GSM.Range("M83").Value = -7420.5
$ws_GSM.Range("N83").Value = -23930.428

# GSM!row132 (G=44008)
$ws_GSM.Range("H132").Value = 4763952
$ws_GSM.Range("I132").Value = 6412035
$ws_GSM.Range("J132").Value = 2821.889
$ws_GSM.Range("K132").Value = 19236105
$ws_GSM.Range("L132").Value = 8465.667000000001
$ws_GSM.Range("M132").Value = -19233575
$ws_GSM.Range("N132").Value = -13525.667

# LTW!row68 (G=12563)
$ws_LTW.Range("H68").Value = 166668460
$ws_LTW.Range("I68").Value = 2225
$ws_LTW.Range("J68").Value = 500000960
$ws_LTW.Range("K68").Value = 2225
$ws_LTW.Range("L68").Value = 500000960
$ws_LTW.Range("M68").Value = -1476
$ws_LTW.Range("N68").Value = -500002458

# LTW!row71 (G=12563)
$ws_LTW.Range("H71").Value = 166668460
$ws_LTW.Range("I71").Value = 2225
$ws_LTW.Range("J71").Value = 500000960
$ws_LTW.Range("K71").Value = 11125
$ws_LTW.Range("L71").Value = 2500004800
$ws_LTW.Range("M71").Value = -7381
$ws_LTW.Range("N71").Value = -2500012288

# LTW!row87 (G=10926)
$ws_LTW.Range("H87").Value = 10525.667
$ws_LTW.Range("I87").Value = 10525.667
$ws_LTW.Range("K87").Value = 10525.667
$ws_LTW.Range("M87").Value = -9402.666999999999

# LTW!row90 (G=10926)
$ws_LTW.Range("H90").Value = 10525.667
$ws_LTW.Range("I90").Value = 10525.667
$ws_LTW.Range("K90").Value = 31577.001
$ws_LTW.Range("M90").Value = -25961.001

# WVR!row96 (G=19977)
$ws_WVR.Range("H96").Value = 2879.5833
$ws_WVR.Range("I96").Value = 2678.111
$ws_WVR.Range("J96").Value = 3484
$ws_WVR.Range("K96").Value = 2678.111
$ws_WVR.Range("L96").Value = 3484
$ws_WVR.Range("M96").Value = -1305.111
$ws_WVR.Range("N96").Value = -6230

# WVR!row121 (G=26316)
$ws_WVR.Range("H121").Value = 30000
$ws_WVR.Range("J121").Value = 30000
$ws_WVR.Range("L121").Value = 30000
$ws_WVR.Range("N121").Value = -33494

# WVR!row122 (G=36208)
$ws_WVR.Range("H122").Value = 1829.8667
$ws_WVR.Range("I122").Value = 1875.7273
$ws_WVR.Range("J122").Value = 1703.75
$ws_WVR.Range("K122").Value = 5627.1819
$ws_WVR.Range("L122").Value = 5111.25
$ws_WVR.Range("M122").Value = -3177.1819
$ws_WVR.Range("N122").Value = -10011.25
